# chore: update Sheets via scheduled runner
# Refreshes cached market-board price / leve-profit figures across several
# item rows on the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets (plain data values,
# no formulas are used in this workbook).

$wb = $excel.ActiveWorkbook

# ALC!15 - Morning Glass of Ether / Ether
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3775.5908
$ws.Range("I15").Value = 3775.5908
$ws.Range("K15").Value = 11326.7724
$ws.Range("M15").Value = -11157.7724

# ALC!17 - One for the Road / Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1007.9773
$ws.Range("J17").Value = 1017.6585
$ws.Range("L17").Value = 3052.9755
$ws.Range("N17").Value = -3388.9755

# ALC!19 - Unbreak My Heart / Roof Tile
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2531.261
$ws.Range("I19").Value = 2418.647
$ws.Range("J19").Value = 2850.3333
$ws.Range("K19").Value = 2418.647
$ws.Range("L19").Value = 2850.3333
$ws.Range("M19").Value = -2243.647
$ws.Range("N19").Value = -3200.3333

# ALC!40 - Stuck in the Moment / Horn Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2679.6667
$ws.Range("I40").Value = 2550.5
$ws.Range("J40").Value = 2699.5386
$ws.Range("K40").Value = 2550.5
$ws.Range("L40").Value = 2699.5386
$ws.Range("M40").Value = -2375.5
$ws.Range("N40").Value = -3049.5386

# ALC!58 - A Matter of Vital Importance / Mega-Potion of Vitality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 521.9286
$ws.Range("I58").Value = 446.69232
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 1340.07696
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -1190.07696
$ws.Range("N58").Value = -4800

# ALC!82 - Rolling on Initiative / Draconian Potion of Dexterity
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = ""  # cell removed (no longer applicable)

# ALC!85 - Darkly Dreaming Dexterity (L) / Draconian Potion of Dexterity
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = ""  # cell removed (no longer applicable)

# ALC!96 - Scroll Down / Grade 1 Reisui of Intelligence
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1854.3529
$ws.Range("I96").Value = 865.2308
$ws.Range("J96").Value = 5069
$ws.Range("K96").Value = 2595.6924
$ws.Range("L96").Value = 15207
$ws.Range("M96").Value = -1222.6924
$ws.Range("N96").Value = -17953

# ALC!132 - Fast-forwarding Flora / Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 183645.61
$ws.Range("I132").Value = 197282.83
$ws.Range("J132").Value = 19999
$ws.Range("K132").Value = 591848.49
$ws.Range("L132").Value = 59997
$ws.Range("M132").Value = -589318.49
$ws.Range("N132").Value = -65057

# ALC!137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I137").Value = 6099.6665
$ws.Range("K137").Value = 18298.9995
$ws.Range("M137").Value = -15748.9995

# ARM!32 - Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6175972.5
$ws.Range("I32").Value = 6175972.5
$ws.Range("K32").Value = 6175972.5
$ws.Range("M32").Value = -6175685.5

# ARM!61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2390477.2
$ws.Range("I61").Value = 2573824
$ws.Range("J61").Value = 6969
$ws.Range("K61").Value = 2573824
$ws.Range("L61").Value = 6969
$ws.Range("M61").Value = -2573612
$ws.Range("N61").Value = -7393

# ARM!132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1188507.5
$ws.Range("I132").Value = 1250634.4
$ws.Range("J132").Value = 8098
$ws.Range("K132").Value = 3751903.2
$ws.Range("L132").Value = 24294
$ws.Range("M132").Value = -3749373.2
$ws.Range("N132").Value = -29354

# ARM!136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2390477.2
$ws.Range("I136").Value = 2573824
$ws.Range("J136").Value = 6969
$ws.Range("K136").Value = 7721472
$ws.Range("L136").Value = 20907
$ws.Range("M136").Value = -7718922
$ws.Range("N136").Value = -26007

# CRP!10 - Spears and Sorcery / Maple Crook
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 316.1
$ws.Range("I10").Value = 244.85715
$ws.Range("K10").Value = 244.85715
$ws.Range("M10").Value = -105.85715

# CRP!11 - Leaving without Leave / Bronze Spear
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2483.5
$ws.Range("J11").Value = 2680.2
$ws.Range("L11").Value = 2680.2
$ws.Range("N11").Value = -2960.2

# CRP!31 - Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 160034.6
$ws.Range("I31").Value = 261506.17
$ws.Range("K31").Value = 261506.17
$ws.Range("M31").Value = -261211.17

# CRP!34 - Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 160034.6
$ws.Range("I34").Value = 261506.17
$ws.Range("K34").Value = 261506.17
$ws.Range("M34").Value = -261304.17

# CRP!58 - You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1238556.2
$ws.Range("I58").Value = 4116448
$ws.Range("J58").Value = 5174
$ws.Range("K58").Value = 4116448
$ws.Range("L58").Value = 5174
$ws.Range("M58").Value = -4116245
$ws.Range("N58").Value = -5580

# CRP!62 - Splinter in the Sewers / Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4312.143
$ws.Range("I62").Value = 2196.25
$ws.Range("J62").Value = 7133.3335
$ws.Range("K62").Value = 2196.25
$ws.Range("L62").Value = 7133.3335
$ws.Range("M62").Value = -1572.25
$ws.Range("N62").Value = -8381.333500000001

# CRP!65 - The Lumber of Their Discontent (L) / Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4312.143
$ws.Range("I65").Value = 2196.25
$ws.Range("J65").Value = 7133.3335
$ws.Range("K65").Value = 10981.25
$ws.Range("L65").Value = 35666.6675
$ws.Range("M65").Value = -7861.25
$ws.Range("N65").Value = -41906.6675

# CRP!132 - Hull Lotta Damage / Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 52910420
$ws.Range("I132").Value = 100013144
$ws.Range("K132").Value = 300039432
$ws.Range("M132").Value = -300036902

# CRP!134 - Wood You Be Quiet / Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 13425.909
$ws.Range("I134").Value = 15996.667
$ws.Range("J134").Value = 1857.5
$ws.Range("K134").Value = 47990.001
$ws.Range("L134").Value = 5572.5
$ws.Range("M134").Value = -45455.001
$ws.Range("N134").Value = -10642.5

# CRP!136 - Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1238556.2
$ws.Range("I136").Value = 4116448
$ws.Range("J136").Value = 5174
$ws.Range("K136").Value = 12349344
$ws.Range("L136").Value = 15522
$ws.Range("M136").Value = -12346794
$ws.Range("N136").Value = -20622

# CRP!140 - Spear Pressure / Claro Walnut Spear
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 89528.5
$ws.Range("J140").Value = 89528.5
$ws.Range("L140").Value = 89528.5
$ws.Range("N140").Value = -99888.5

# CUL!32 - Convalescence Precedes Essence / Ginger Cookie
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 500000700
$ws.Range("J32").Value = 1000000000
$ws.Range("L32").Value = 3000000000
$ws.Range("N32").Value = -3000000566

# CUL!75 - Breakfast of Champions / Emerald Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7363.0835
$ws.Range("I75").Value = 3000
$ws.Range("J75").Value = 7759.727
$ws.Range("K75").Value = 9000
$ws.Range("L75").Value = 23279.181
$ws.Range("M75").Value = -8002
$ws.Range("N75").Value = -25275.181

# CUL!78 - Emerald Soup for the Soul (L) / Emerald Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 7363.0835
$ws.Range("I78").Value = 3000
$ws.Range("J78").Value = 7759.727
$ws.Range("K78").Value = 27000
$ws.Range("L78").Value = 69837.54300000001
$ws.Range("M78").Value = -22008
$ws.Range("N78").Value = -79821.54300000001

# CUL!80 - Saucy for a Suitor / Hollandaise Sauce
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3600

# CUL!83 - Saved by the Sauce (L) / Hollandaise Sauce
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3600

# GSM!18 - Gorgeous Gorget / Brass Gorget
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4050.7
$ws.Range("I18").Value = 2042.8334
$ws.Range("J18").Value = 7062.5
$ws.Range("K18").Value = 2042.8334
$ws.Range("L18").Value = 7062.5
$ws.Range("M18").Value = -1749.8334
$ws.Range("N18").Value = -7648.5

# GSM!132 - On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19102700
$ws.Range("I132").Value = 28921112
$ws.Range("J132").Value = 11346.223
$ws.Range("K132").Value = 86763336
$ws.Range("L132").Value = 34038.669
$ws.Range("M132").Value = -86760806
$ws.Range("N132").Value = -39098.669

# LTW!11 - A Thorn in One's Hide / Leather Mitts
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 53996.668
$ws.Range("J11").Value = 53996.668
$ws.Range("L11").Value = 53996.668
$ws.Range("N11").Value = -54276.668

# LTW!46 - Supply Side Logic / Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1636
$ws.Range("J46").Value = 2028.2
$ws.Range("L46").Value = 2028.2
$ws.Range("N46").Value = -2404.2

# LTW!93 - Hide to Go Seek / Gagana Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 290.66666
$ws.Range("I93").Value = 201.5
$ws.Range("J93").Value = 469
$ws.Range("K93").Value = 201.5
$ws.Range("L93").Value = 469
$ws.Range("M93").Value = 1046.5
$ws.Range("N93").Value = -2965

# LTW!132 - Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1661884.6
$ws.Range("I132").Value = 2488541.8
$ws.Range("K132").Value = 7465625.399999999
$ws.Range("M132").Value = -7463095.399999999

# WVR!18 - Welcome to the Cotton Club / Cotton Halfgloves
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 34000
$ws.Range("J18").Value = 34000
$ws.Range("L18").Value = 34000
$ws.Range("N18").Value = -34346

# WVR!132 - Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3303973.5
$ws.Range("I132").Value = 4285467.5
$ws.Range("J132").Value = 8958.214
$ws.Range("K132").Value = 12856402.5
$ws.Range("L132").Value = 26874.642
$ws.Range("M132").Value = -12853872.5
$ws.Range("N132").Value = -31934.642
